$d = $word.ActiveDocument
$dash = [char]0x2013

# 1) Current address: expand the street/area details and fix the pin code.
$oldAddr = "Curr Add : 16, Keshav Chaya Apt, Bhau Patil Road, Khadki, Pune " + $dash + " 411003"
$newAddr = "Curr Add : R-25, Dhruv Darshan Soc., Near PCCOE, Akurdi, Pune " + $dash + " 411044"
$d.Content.Find.Execute($oldAddr, $true, $false, $false, $false, $false, $true, 1, $false, $newAddr, 2)

# 2) CGPA value update.
$d.Content.Find.Execute("CGPA : 7.2", $true, $false, $false, $false, $false, $true, 1, $false, "CGPA : 9.03", 2)

# 3) Project title rename.
$d.Content.Find.Execute("Parking Automation", $true, $false, $false, $false, $false, $true, 1, $false, "Thieves Recognition", 2)

# 4) Project description: append extra technologies used.
$d.Content.Find.Execute("Python, OpenCV, Deep Learning", $true, $false, $false, $false, $false, $true, 1, $false, "Python, OpenCV, Deep Learning, Django, MySQL", 2)

# 5) Certification entry replaced.
$d.Content.Find.Execute("Intro of ML (NPTEL)", $true, $false, $false, $false, $false, $true, 1, $false, "Hacker Rank SQL Advance (Hacker Rank)", 2)
